$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "8:15-8:20"
$ws.Range("C3").Value = "8:20-8:25"
$ws.Range("C6").Value = "22:40-22:45"
$ws.Range("C7").Value = "22:45-22:50"

$ws.Range("C8").Select()
